$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be written as text (not auto-converted to a
    # number/date by Excel) while keeping its original (unstyled) look:
    # apply a text number format only long enough to assign the value,
    # then restore the cell to the default "Normal" style so no stray
    # formatting is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "58.385.84"
Set-TextValue $ws.Range("E2") "  -0.95%  "

Set-TextValue $ws.Range("D3") "2.480.02"
Set-TextValue $ws.Range("E3") "  -0.67%  "

Set-TextValue $ws.Range("E4") "  +0.08%  "

Set-TextValue $ws.Range("D5") "525.62"
Set-TextValue $ws.Range("E5") "  -2.08%  "

Set-TextValue $ws.Range("D6") "133.10"
Set-TextValue $ws.Range("E6") "  -3.44%  "

Set-TextValue $ws.Range("E7") "  +0.39%  "

Set-TextValue $ws.Range("D8") "0.560"
Set-TextValue $ws.Range("E8") "  -1.09%  "

Set-TextValue $ws.Range("D9") "0.1000"
Set-TextValue $ws.Range("E9") "  -1.08%  "

Set-TextValue $ws.Range("D10") "0.157"
Set-TextValue $ws.Range("E10") "  -1.98%  "

Set-TextValue $ws.Range("D11") "5.36"
Set-TextValue $ws.Range("E11") "  +0.19%  "

Set-TextValue $ws.Range("D12") "0.340"
Set-TextValue $ws.Range("E12") "  -1.75%  "

Set-TextValue $ws.Range("D13") "2.921.38"
Set-TextValue $ws.Range("E13") "  -0.87%  "

Set-TextValue $ws.Range("D14") "58.398.14"
Set-TextValue $ws.Range("E14") "  -0.76%  "

Set-TextValue $ws.Range("D15") "22.29"
Set-TextValue $ws.Range("E15") "  -3.76%  "

Set-TextValue $ws.Range("D16") "0.0000136"
Set-TextValue $ws.Range("E16") "  -2.32%  "

Set-TextValue $ws.Range("D17") "2.482.92"
Set-TextValue $ws.Range("E17") "  -1.26%  "

Set-TextValue $ws.Range("D18") "10.86"
Set-TextValue $ws.Range("E18") "  -2.03%  "

Set-TextValue $ws.Range("D19") "4.20"
Set-TextValue $ws.Range("E19") "  -2.03%  "

Set-TextValue $ws.Range("D20") "320.08"
Set-TextValue $ws.Range("E20") "  -1.65%  "

Set-TextValue $ws.Range("E21") "  +0.09%  "

Set-TextValue $ws.Range("D22") "5.79"
Set-TextValue $ws.Range("E22") "  -1.41%  "

Set-TextValue $ws.Range("D23") "63.96"
Set-TextValue $ws.Range("E23") "  -2.38%  "

Set-TextValue $ws.Range("D24") "0.412"
Set-TextValue $ws.Range("E24") "  -2.35%  "

Set-TextValue $ws.Range("B25") "Binance-PegBSC-USD"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D25") "1.00"
Set-TextValue $ws.Range("E25") "  +0.24%  "

Set-TextValue $ws.Range("B26") "Kaspa"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D26") "0.162"
Set-TextValue $ws.Range("E26") "  -2.10%  "

Set-TextValue $ws.Range("D27") "7.42"
Set-TextValue $ws.Range("E27") "  -2.98%  "

Set-TextValue $ws.Range("D28") "0.0₃0748"
Set-TextValue $ws.Range("E28") "  -3.78%  "

Set-TextValue $ws.Range("D29") "6.40"
Set-TextValue $ws.Range("E29") "  -4.22%  "

Set-TextValue $ws.Range("D30") "1.71"
Set-TextValue $ws.Range("E30") "  -3.62%  "

Set-TextValue $ws.Range("D31") "166.57"
Set-TextValue $ws.Range("E31") "  -1.42%  "

Set-TextValue $ws.Range("B32") "Fetch.AI"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D32") "1.13"
Set-TextValue $ws.Range("E32") "  -6.09%  "

Set-TextValue $ws.Range("B33") "USDe"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D33") "0.999"
Set-TextValue $ws.Range("E33") "  +0.11%  "

Set-TextValue $ws.Range("D34") "1.00"
Set-TextValue $ws.Range("E34") "  +0.54%  "

Set-TextValue $ws.Range("D35") "18.20"
Set-TextValue $ws.Range("E35") "  -1.95%  "

Set-TextValue $ws.Range("D36") "1.35"
Set-TextValue $ws.Range("E36") "  -9.69%  "

Set-TextValue $ws.Range("D37") "3.96"
Set-TextValue $ws.Range("E37") "  -3.64%  "

Set-TextValue $ws.Range("D38") "1.49"
Set-TextValue $ws.Range("E38") "  -4.66%  "

Set-TextValue $ws.Range("D39") "3.52"
Set-TextValue $ws.Range("E39") "  -3.27%  "

Set-TextValue $ws.Range("D40") "0.787"
Set-TextValue $ws.Range("E40") "  -4.51%  "

Set-TextValue $ws.Range("D41") "275.43"
Set-TextValue $ws.Range("E41") "  -3.56%  "

Set-TextValue $ws.Range("B42") "RenderToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D42") "4.93"
Set-TextValue $ws.Range("E42") "  -8.14%  "

Set-TextValue $ws.Range("B43") "Mantle"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D43") "0.596"
Set-TextValue $ws.Range("E43") "  -1.27%  "

Set-TextValue $ws.Range("B44") "Aave"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D44") "127.38"
Set-TextValue $ws.Range("E44") "  -2.92%  "

Set-TextValue $ws.Range("B45") "Stellar"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D45") "0.0913"
Set-TextValue $ws.Range("E45") "  -2.09%  "

Set-TextValue $ws.Range("B46") "Hedera"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D46") "0.0492"
Set-TextValue $ws.Range("E46") "  -3.98%  "

Set-TextValue $ws.Range("B47") "VeChain"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D47") "0.0215"
Set-TextValue $ws.Range("E47") "  -3.00%  "

Set-TextValue $ws.Range("B48") "InjectiveProtocol"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D48") "17.00"
Set-TextValue $ws.Range("E48") "  -3.13%  "

Set-TextValue $ws.Range("B49") "Maker"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D49") "1.736.07"
Set-TextValue $ws.Range("E49") "  -1.64%  "

Set-TextValue $ws.Range("B50") "BitgetToken"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
Set-TextValue $ws.Range("D50") "0.978"
Set-TextValue $ws.Range("E50") "  -1.75%  "

Set-TextValue $ws.Range("B51") "ZEEBU"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
Set-TextValue $ws.Range("D51") "4.73"
Set-TextValue $ws.Range("E51") "  -0.88%  "

